$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows starting at row 18; this pushes the existing
# "Sector Distribution Details" block (old rows 18-33) down to rows 24-39.
$ws.Range("A18:A23").EntireRow.Insert()

# --- New "Size Distribution Details" table (rows 17-21) ---

# Header row (bold, like the other section headers)
$ws.Range("B17").Value = "Number of employees"
$ws.Range("B17").Font.Bold = $true
$ws.Range("C17").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("C17").Font.Bold = $true
$ws.Range("D17").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("D17").Font.Bold = $true

# Micro row
$ws.Range("A18").Value = "Micro"
$ws.Range("B18").Font.Bold = $false
$ws.Range("C18").Font.Bold = $false
$ws.Range("D18").Font.Bold = $false

# Small row
$ws.Range("A19").Value = "Small"
$ws.Range("B19").Value = "5-10"
$ws.Range("C19").Font.Bold = $false
$ws.Range("D19").Font.Bold = $false

# Medium row
$ws.Range("A20").Value = "Medium"
$ws.Range("B20").Font.Bold = $false
$ws.Range("C20").Font.Bold = $false
$ws.Range("D20").Font.Bold = $false

# Large row
$ws.Range("A21").Value = "Large"
$ws.Range("B21").Font.Bold = $false
$ws.Range("C21").Font.Bold = $false
$ws.Range("D21").Font.Bold = $false

# --- Fix up the hyperlink, whose anchor cell moved from A28 to A34 ---
$ws.Range("A28").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A34"), "http://archives.entreprendre-mediterranee.com/documents/colloque-Acim-tunis-2006/Session4/Yeye-Burkina-Faso.pdf")
